$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Map of cell address -> new value (text, matching the original "numeric-looking
# text" storage used throughout this sheet). Values are forced to stay text by
# temporarily switching the cell to a Text number format before assignment,
# then restoring the "Normal" style so the cell's style index is unchanged.
$changes = @{
    "B11" = "41.79"
    "C11" = "8.35"
    "D11" = "50.15"
    "B33" = "34.23"
    "C33" = "3.96"
    "D33" = "38.19"
    "B34" = "21.43"
    "C34" = "44.99"
    "D34" = "66.43"
    "B36" = "89.36"
    "C36" = "10.34"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    $cell.Style = "Normal"
}
